# Penambahan API dan Update Data API Catalogue
# Insert a new "float" restriction row (row 7) into the JSON Schema
# Restriction sheet, pushing the existing "date" and "timestamptz" rows
# down by one, and refresh the frozen-pane / selection view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- keep the sheet's used-range anchored at row 1 -------------------
# Row 1 only carries a custom (very small) row height and no cell
# content. Touching it briefly makes the engine keep it as part of the
# worksheet dimension once new rows are inserted below, matching the
# original top-left anchor of the sheet.
$ws.Range("B1").Value = 1
$ws.Range("B1").Value = $null
$ws.Rows("1:1").RowHeight = 7.5

# --- insert the new row for the "float" type --------------------------
$ws.Rows("7:7").Insert()

# Copy the formatting (styles/borders) of the row above into the new
# row so the new cells line up with the rest of the table.
$ws.Range("B6:E6").Copy()
$ws.Range("B7:E7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B7").Value = "float"
$ws.Range("C7").Value = """type"": ""number"""

# --- refresh frozen pane / selection state -----------------------------
$ws.Application.ActiveWindow.Panes.Item(4).ScrollColumn = 3
$ws.Range("C7").Select()
